# Add calc of 2020 ART/VMMC to Excel sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two new columns before column I (old I.. shifts to K.., old L shifts to N)
$ws.Range("I1:J1").EntireColumn.Insert()

# 2. New header values in H2:J2
$ws.Range("H2").Value = 2020
$ws.Range("I2").Value = 2030
$ws.Range("J2").Value = "slope"

# 3. New "interim calculation..." note in I7 (merged I7:I12)
$ws.Range("I7").Value = "interim calculation assuming linear scale-up to 90-90-90 by 2030"

# 4. Target (90-90-90) proportions for each age band
$ws.Range("I3").Value = 0.7
$ws.Range("I4").Value = 0.7
$ws.Range("I5").Value = 0.7
$ws.Range("I6").Value = 0.7

# 5. slope formulas (J3:J6) -- entered individually (not shared)
$ws.Range("J3").Formula = "=(I3-G3)/((I`$2-1)-(G`$2-1))"
$ws.Range("J4").Formula = "=(I4-G4)/((I`$2-1)-(G`$2-1))"
$ws.Range("J5").Formula = "=(I5-G5)/((I`$2-1)-(G`$2-1))"
$ws.Range("J6").Formula = "=(I6-G6)/((I`$2-1)-(G`$2-1))"

# 6. 2020 interpolation formulas (H3 standalone, H4:H6 filled as a shared group)
$ws.Range("H3").Formula = "=G3+(H`$2-G`$2)*J3"
$ws.Range("H4:H6").Formula = "=G4+(H`$2-G`$2)*J4"

# 7. Number format for H3:H6 (six decimal places)
$ws.Range("H3:H6").NumberFormat = "0.000000"

# 8. Font color (red) for the new calculation cells
$ws.Range("H2:J6").Font.Color = 255

# 9. Merge I7:I12 and format (left align + wrap text)
$ws.Range("I7:I12").Merge()
$ws.Range("I7:I12").HorizontalAlignment = -4131
$ws.Range("I7:I12").WrapText = $true
$ws.Range("I7:I12").Font.Color = 255

# 10. Wrap text for I13:I14 (continuation of merged note column formatting)
$ws.Range("I13:I14").WrapText = $true

# 11. Column I width
$ws.Range("I1").ColumnWidth = 11.1640625

# 12. Row 7 explicit height (matches new default row height)
$ws.Range("A7:N7").RowHeight = 15.75

# 13. Page setup or orientation
$ws.PageSetup.Orientation = 1

# 14. Update selection to match the authored file
$ws.Range("M28").Select()
